$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Add double line spacing (w:line="480" w:lineRule="auto") to the
#    three body paragraphs of the reflection.
# ---------------------------------------------------------------------
$bodyParas = @(
  "In revising, I did not make major structural changes",
  "Above all else, the strength of this essay",
  "I recognize that while the portion of the essay"
)
foreach ($needle in $bodyParas) {
  $rng = $d.Content
  $found = $rng.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
  if ($found) {
    $para = $rng.Paragraphs(1)
    $para.LineSpacingRule = 5   # wdLineSpaceMultiple
    $para.LineSpacing = 24      # 24pt -> stored as w:line="480" (twentieths of a point)
  }
}

# ---------------------------------------------------------------------
# 2) Paragraph "In revising...": split the big run right after
#    "...appeal to a" and drop a (re-seated) _GoBack bookmark there.
#    Word only keeps a single _GoBack bookmark, so this automatically
#    removes the one that used to live in the next paragraph.
# ---------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("making the story appeal to a", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitPoint = $d.Range($rng.End, $rng.End)
$d.Bookmarks.Add("_GoBack", $splitPoint)

# ---------------------------------------------------------------------
# 3) Paragraph "Above all else...": the old bookmark split the run into
#    two; now that the bookmark has moved away, re-merge that text back
#    into a single run via a no-op Find/Replace across both halves.
# ---------------------------------------------------------------------
$rng = $d.Content
$mergeText = "catered more specially toward one of them."
$rng.Find.Execute($mergeText, $true, $false, $false, $false, $false, $true, 1, $false, $mergeText, 2) | Out-Null

# ---------------------------------------------------------------------
# 4) Last paragraph: insert a <w:lastRenderedPageBreak/> right before
#    "grammatical feedback...". Build the split by temporarily breaking
#    the paragraph in two, splicing the page-break-marked run in via
#    InsertXML, rejoining the paragraph, then removing the now-duplicate
#    trailing copy of the text.
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("very much welcome ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitPoint = $d.Range($rng.End, $rng.End)
$splitPoint.InsertParagraphAfter()

$tailText = "grammatical feedback, since there" + [char]0x2019 + "s no way I" + [char]0x2019 + "ll improve it if it reads fine to me and nothing more comes of it."

$rng = $d.Content
$rng.Find.Execute("grammatical feedback, since there", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$insertPt = $d.Range($rng.Start, $rng.Start)

$runProps = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/></w:rPr>'
$escapedTail = $tailText -replace '&','&amp;' -replace '<','&lt;' -replace '>','&gt;'
$xml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' +
       '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
       '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' +
       '<w:r>' + $runProps + '<w:lastRenderedPageBreak/><w:t xml:space="preserve">' + $escapedTail + '</w:t></w:r>' +
       '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertPt.InsertXML($xml)

# Rejoin the temporary paragraph split (delete the paragraph mark).
$rng = $d.Content
$rng.Find.Execute("very much welcome ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$pilcrow = $d.Range($rng.End, $rng.End + 1)
$pilcrow.Delete()

# Remove the now-duplicated original copy of the tail text (2nd occurrence).
$rng = $d.Content
$rng.Find.Execute($tailText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$afterFirst = $d.Range($rng.End, $d.Content.End)
$afterFirst.Find.Execute($tailText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$afterFirst.Delete()

Write-Host "done"
